$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.077228333333333
$ws.Range("H2").Value = 3.231685
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.02564166666666666
$ws.Range("N2").Value = 0.07692499999999999
$ws.Range("O2").Value = 0.0006780701807970013
$ws.Range("P2").Value = 0.0006780701807970013
$ws.Range("Q2").Value = 0.02762192984722222
$ws.Range("R2").Value = 0.248597368625
$ws.Range("S2").Value = 0.0006780701807970013
$ws.Range("T2").Value = 0.0006780701807970013

# Row 3
$ws.Range("G3").Value = 1.077228333333333
$ws.Range("H3").Value = 3.231685
$ws.Range("M3").Value = 0.01112833333333333
$ws.Range("O3").Value = 0.0002942784918545062
$ws.Range("P3").Value = 0.0002942784918545062
$ws.Range("Q3").Value = 0.01198775596944444
$ws.Range("R3").Value = 0.107889803725
$ws.Range("S3").Value = 0.0002942784918545062
$ws.Range("T3").Value = 0.0002942784918545062

# Row 4
$ws.Range("G4").Value = 1.077228333333333
$ws.Range("H4").Value = 3.231685
$ws.Range("O4").Value = 0.0009221225577320236
$ws.Range("P4").Value = 0.0009221225577320235
$ws.Range("Q4").Value = 0.03756367013555555
$ws.Range("R4").Value = 0.3380730312199999
$ws.Range("S4").Value = 0.0009221225577320236
$ws.Range("T4").Value = 0.0009221225577320235

# Row 5
$ws.Range("G5").Value = 1.077228333333333
$ws.Range("H5").Value = 3.231685
$ws.Range("M5").Value = 37.74401233333333
$ws.Range("N5").Value = 113.232037
$ws.Range("O5").Value = 0.9981055287696164
$ws.Range("P5").Value = 0.9981055287696164
$ws.Range("Q5").Value = 40.65891949914944
$ws.Range("R5").Value = 365.9302754923449
$ws.Range("S5").Value = 0.9981055287696164
$ws.Range("T5").Value = 0.9981055287696164

Write-Output "applied updates"
